$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RCL")

# Row 6: Change in inventories
$ws.Range("B6").Value = 38967000.0
$ws.Range("C6").Value = 27077000.0
$ws.Range("D6").Value = 24782000.0
$ws.Range("E6").Value = 15414000.0
$ws.Range("F6").Value = -19765000.0

# Row 8: Change in payables and accrued liability
$ws.Range("B8").Value = 1776000000.0
$ws.Range("C8").Value = 2849000000.0
$ws.Range("D8").Value = 2456993000.0
$ws.Range("E8").Value = 2132211000.0
$ws.Range("F8").Value = 1529137000.0

# Row 29: Capital Stock Change
$ws.Range("B29").Value = 3048609000.0
$ws.Range("C29").Value = 1431759000.0
$ws.Range("D29").Value = -100264000.0
$ws.Range("E29").Value = -99077000.0

$wb.Save()
